$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.652.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.094.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5149"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4368"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.43"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09161"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.165"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.48"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.090.37"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.745"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.183"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.48"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.93"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06673"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.205"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.708.79"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.339.02"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.87"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.50"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.488"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.44"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.125"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1049"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.660"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.197"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.941"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.299"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.94%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02570"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06662"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6948"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.45"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.327"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2215"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6810"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.32"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.300"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.612"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000351"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.216"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.214"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.12"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.68%  "
